# Adds reference to project organization doc.
#
# Inserts a new list-paragraph (same numbered-list level as its sibling
# bullets: ilvl=1, numId=4) right after the "ProcessingLog.xlsx" bullet
# and right before the "Raw Data QC" bullet. The new paragraph has two
# runs: a plain-text lead-in and a bold filename.

$d = $word.ActiveDocument

# Locate the "Raw Data QC" paragraph - the new paragraph must be inserted
# immediately before it (i.e. immediately after the ProcessingLog.xlsx
# paragraph that currently precedes it).
$rng = $d.Content
$found = $rng.Find.Execute("Raw Data QC", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$targetPara = $rng.Paragraphs(1)
$targetIndex = $targetPara.Index

# Insert a new (empty) paragraph right before it. This new paragraph
# inherits formatting from the "Raw Data QC" paragraph, which is
# conveniently un-bold already.
$targetPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($targetIndex)
$insertRange = $newPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>' + `
    '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' + `
    '<w:t xml:space="preserve">For more on file organization, see </w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' + `
    '<w:t>03_Project_Organization.docx</w:t>' + `
  '</w:r>' + `
'</w:p>' + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$null = $insertRange.InsertXML($xml)
